# Update attendance detail cells from 0 to 1 per row, per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> list of columns to flip from 0 to 1
$updates = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("D", "E")
    6  = @("H")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("H")
    11 = @("H")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
